# Weekly fruit/vegetable price update: two new daily price records were
# added at the top of the Coliflor (Vega Monumental Concepción) table.
# This pushes all existing records down by two rows (old row 285 -> new
# row 287, ... old row 387 -> new row 389), and the two freshly reported
# rows are written into the newly opened 285:286 slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 285, shifting all
# rows from 285 downward (including the former last row, 387) down by 2.
$ws.Rows("285:286").Insert()

# --- New row 285 -------------------------------------------------------
$ws.Range("A285").Value = 11
$ws.Range("B285").Value = "Vega Monumental Concepción"
$ws.Range("C285").Value = "Bíobío"
$ws.Range("D285").Value = 44992
$ws.Range("E285").Value = 8
$ws.Range("F285").Value = 100112008
$ws.Range("G285").Value = "Coliflor"
$ws.Range("H285").Value = "Sin especificar"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 2000
$ws.Range("K285").Value = 1100
$ws.Range("L285").Value = 1200
$ws.Range("M285").Value = 1150
$ws.Range("N285").Value = "$/unidad"
$ws.Range("O285").Value = "Región Metropolitana"
$ws.Range("P285").Value = 1150
$ws.Range("Q285").Value = 1
$ws.Range("R285").Value = "Hortaliza"

# --- New row 286 -------------------------------------------------------
$ws.Range("A286").Value = 11
$ws.Range("B286").Value = "Vega Monumental Concepción"
$ws.Range("C286").Value = "Bíobío"
$ws.Range("D286").Value = 44992
$ws.Range("E286").Value = 8
$ws.Range("F286").Value = 100112008
$ws.Range("G286").Value = "Coliflor"
$ws.Range("H286").Value = "Sin especificar"
$ws.Range("I286").Value = "Segunda"
$ws.Range("J286").Value = 1000
$ws.Range("K286").Value = 900
$ws.Range("L286").Value = 900
$ws.Range("M286").Value = 900
$ws.Range("N286").Value = "$/unidad"
$ws.Range("O286").Value = "Región Metropolitana"
$ws.Range("P286").Value = 900
$ws.Range("Q286").Value = 1
$ws.Range("R286").Value = "Hortaliza"
